$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (order chosen to match shared-string insertion order)
# A2/E2 are prefixed with a leading apostrophe so Excel keeps storing them as
# literal text (quote-prefixed) instead of re-interpreting them as a number,
# matching the original cell formatting (s="1" / s="3").
$ws.Range("C2").Value = "su"
$ws.Range("F2").Value = "Inclusión de Riesgo"
$ws.Range("D2").Value = "silverarrow"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("A2").Value = "'i-preproducciongestion.segurossura.com.ar"
$ws.Range("E2").Value = "'04104013552"

# Remove the hyperlink on B2 (keep text/value, drop the link)
$ws.Range("B2").Hyperlinks.Delete()

# Update the selection to E5
$ws.Range("E5").Select()
